# email service added: remove the "Amar Jadhav / Senior Software Engineer"
# appreciation entry (row 2) from the Appreciations sheet; remaining rows
# shift up automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Appreciations")
$ws.Rows.Item(2).Delete()
